# Generate Report for Handback
# Update the localization-status report to reflect that the handback for
# c6c61c5a-1579-4e80-85c7-39e8865daba8.md is now complete ("Handed back:
# in sync with en-US") instead of "Ready for handoff", refresh the Latest
# Handback DateTime for each locale, and clear the stale Error Detail.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for c6c61c5a...md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row for c6c61c5a...md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-28 02:48:23"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row for c6c61c5a...md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-28 02:48:29"
$wsDeDe.Range("P3").Value = ""
